# Update Poker - Person Info worksheet with refreshed query values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 - Mark
$ws.Range("D3").Value = 152
$ws.Range("I3").Value = 547
$ws.Range("J3").Value = 3.6

# Row 5 - Pepe
$ws.Range("D5").Value = 111
$ws.Range("I5").Value = 409
$ws.Range("J5").Value = 3.68

# Row 6 - Prashant
$ws.Range("D6").Value = 48
$ws.Range("I6").Value = 180
$ws.Range("J6").Value = 3.75

# Row 7 - Richard
$ws.Range("D7").Value = 151
$ws.Range("I7").Value = 646

# Row 10 - Andy
$ws.Range("D10").Value = 209
$ws.Range("I10").Value = 869
$ws.Range("J10").Value = 4.16

# Row 11 - Anthony
$ws.Range("D11").Value = 129
$ws.Range("I11").Value = 506
